# "getting sentdex's code to work, and view visualization on tensorboard"
# The "Model Performance" sheet is renamed to reflect the optimized model,
# and the cursor is left on the new Accuracy figures (H7:H8) after
# reviewing them across the workbook's sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AAPL ML Perf")
$ws2 = $wb.Worksheets.Item("SPY ML Perf")
$ws3 = $wb.Worksheets.Item("Model Performance")

$ws3.Name = "Optimized Model Performance"

$ws1.Activate()
$ws1.Range("H7:H8").Select()

$ws2.Activate()
$ws2.Range("H7:H8").Select()

# Leave the renamed, Accuracy-reviewing sheet active/selected last so it
# stays the active tab, matching its tabSelected state.
$ws3.Activate()
$ws3.Range("H7:H8").Select()
